$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 658 (existing rows 658+ shift down by 2)
$ws.Rows.Item(658).Insert()
$ws.Rows.Item(658).Insert()

# New row 658: 2026/01/19 (Mon) 23:00 -> rank 13
$ws.Range("A658").NumberFormat = "@"
$ws.Range("A658").Value = "2026/01/19"
$ws.Range("B658").Value = "月"
$ws.Range("C658").Value = 23
$ws.Range("D658").Value = 13

# New row 659: 2026/01/20 (Tue) 2:00 -> rank 15
$ws.Range("A659").NumberFormat = "@"
$ws.Range("A659").Value = "2026/01/20"
$ws.Range("B659").Value = "火"
$ws.Range("C659").Value = 2
$ws.Range("D659").Value = 15

# Restore default style so no explicit style index is left on the new cells
$ws.Range("A658:A659").Style = "Normal"
